$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H40").Value = 4170811.2
$ws.Range("I40").Value = 11114653
$ws.Range("J40").Value = 4506.067
$ws.Range("K40").Value = 11114653
$ws.Range("L40").Value = 4506.067
$ws.Range("M40").Value = -11114478
$ws.Range("N40").Value = -4856.067
$ws.Range("H58").Value = 260.8889
$ws.Range("I58").Value = 260.8889
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 782.6667
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -632.6667
$ws.Range("N58").ClearContents()
$ws.Range("H64").Value = 5539
$ws.Range("I64").Value = 5048.75
$ws.Range("K64").Value = 5048.75
$ws.Range("M64").Value = -4800.75
$ws.Range("H67").Value = 5539
$ws.Range("I67").Value = 5048.75
$ws.Range("K67").Value = 5048.75
$ws.Range("M67").Value = -4190.75
$ws.Range("H76").Value = 3341.5386
$ws.Range("I76").Value = 2536.5
$ws.Range("K76").Value = 2536.5
$ws.Range("M76").Value = -2221.5
$ws.Range("H79").Value = 3341.5386
$ws.Range("I79").Value = 2536.5
$ws.Range("K79").Value = 2536.5
$ws.Range("M79").Value = -1444.5
$ws.Range("H80").Value = 4063.7856
$ws.Range("I80").Value = 2954
$ws.Range("K80").Value = 8862
$ws.Range("M80").Value = -7864
$ws.Range("H83").Value = 4063.7856
$ws.Range("I83").Value = 2954
$ws.Range("K83").Value = 26586
$ws.Range("M83").Value = -21594
$ws.Range("H98").Value = 3670.9
$ws.Range("I98").Value = 1720.8572
$ws.Range("J98").Value = 8221
$ws.Range("K98").Value = 1720.8572
$ws.Range("L98").Value = 8221
$ws.Range("M98").Value = -222.8571999999999
$ws.Range("N98").Value = -11217
$ws.Range("H121").Value = 6250
$ws.Range("J121").Value = 6250
$ws.Range("L121").Value = 18750
$ws.Range("N121").Value = -22244
$ws.Range("H122").Value = 3670.9
$ws.Range("I122").Value = 1720.8572
$ws.Range("J122").Value = 8221
$ws.Range("K122").Value = 5162.571599999999
$ws.Range("L122").Value = 24663
$ws.Range("M122").Value = -2712.571599999999
$ws.Range("N122").Value = -29563
$ws.Range("H135").Value = 3371.9473
$ws.Range("I135").Value = 1926.2
$ws.Range("J135").Value = 4978.3335
$ws.Range("K135").Value = 17335.8
$ws.Range("L135").Value = 44805.0015
$ws.Range("M135").Value = -14800.8
$ws.Range("N135").Value = -49875.0015

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H74").Value = 1403.1
$ws.Range("I74").Value = 824
$ws.Range("J74").Value = 1896.4073
$ws.Range("K74").Value = 824
$ws.Range("L74").Value = 1896.4073
$ws.Range("M74").Value = 50
$ws.Range("N74").Value = -3644.4073
$ws.Range("H77").Value = 1403.1
$ws.Range("I77").Value = 824
$ws.Range("J77").Value = 1896.4073
$ws.Range("K77").Value = 4120
$ws.Range("L77").Value = 9482.0365
$ws.Range("M77").Value = 248
$ws.Range("N77").Value = -18218.0365
$ws.Range("H132").Value = 5827.9473
$ws.Range("I132").Value = 2128.5789
$ws.Range("K132").Value = 6385.736699999999
$ws.Range("M132").Value = -3855.736699999999

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H86").Value = 1962300.5
$ws.Range("I86").Value = 2899846
$ws.Range("K86").Value = 2899846
$ws.Range("M86").Value = -2898723
$ws.Range("H89").Value = 1962300.5
$ws.Range("I89").Value = 2899846
$ws.Range("K89").Value = 14499230
$ws.Range("M89").Value = -14493614
$ws.Range("H94").Value = 2148.5833
$ws.Range("J94").Value = 2677.2
$ws.Range("L94").Value = 2677.2
$ws.Range("N94").Value = -3579.2
$ws.Range("H105").Value = 4156.125
$ws.Range("I105").Value = 3587.25
$ws.Range("K105").Value = 3587.25
$ws.Range("M105").Value = -1840.25
$ws.Range("H134").Value = 1844.8909
$ws.Range("I134").Value = 1276.3158
$ws.Range("K134").Value = 3828.9474
$ws.Range("M134").Value = -1293.9474

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H31").Value = 2130.3076
$ws.Range("I31").Value = 1647.4
$ws.Range("K31").Value = 1647.4
$ws.Range("M31").Value = -1352.4
$ws.Range("H34").Value = 2130.3076
$ws.Range("I34").Value = 1647.4
$ws.Range("K34").Value = 1647.4
$ws.Range("M34").Value = -1445.4
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H132").Value = 7347
$ws.Range("I132").Value = 1003.6667
$ws.Range("K132").Value = 3011.0001
$ws.Range("M132").Value = -481.0001000000002

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 13207687
$ws.Range("I4").Value = 14967742
$ws.Range("K4").Value = 44903226
$ws.Range("M4").Value = -44903114
$ws.Range("H34").Value = 1235.7142
$ws.Range("I34").Value = 1050
$ws.Range("K34").Value = 3150
$ws.Range("M34").Value = -3066
$ws.Range("H122").Value = 360
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H80").Value = 6821.0557
$ws.Range("I80").Value = 8554
$ws.Range("J80").Value = 4097.857
$ws.Range("K80").Value = 8554
$ws.Range("L80").Value = 4097.857
$ws.Range("M80").Value = -7556
$ws.Range("N80").Value = -6093.857
$ws.Range("H83").Value = 6821.0557
$ws.Range("I83").Value = 8554
$ws.Range("J83").Value = 4097.857
$ws.Range("K83").Value = 42770
$ws.Range("L83").Value = 20489.285
$ws.Range("M83").Value = -37778
$ws.Range("N83").Value = -30473.285
$ws.Range("H107").Value = 649.8
$ws.Range("I107").Value = 427.7647
$ws.Range("J107").Value = 1121.625
$ws.Range("K107").Value = 427.7647
$ws.Range("L107").Value = 1121.625
$ws.Range("M107").Value = 1492.2353
$ws.Range("N107").Value = -4961.625
$ws.Range("H122").Value = 4259.9165
$ws.Range("I122").Value = 3456.125
$ws.Range("K122").Value = 10368.375
$ws.Range("M122").Value = -7918.375
$ws.Range("H132").Value = 4177.95
$ws.Range("I132").Value = 3983.0725
$ws.Range("K132").Value = 11949.2175
$ws.Range("M132").Value = -9419.217500000001
$ws.Range("H135").Value = 81099.5
$ws.Range("J135").Value = 81099.5
$ws.Range("L135").Value = 81099.5
$ws.Range("N135").Value = -91239.5

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H46").Value = 11328.25
$ws.Range("I46").Value = 3662
$ws.Range("J46").Value = 22061
$ws.Range("K46").Value = 3662
$ws.Range("L46").Value = 22061
$ws.Range("M46").Value = -3474
$ws.Range("N46").Value = -22437
$ws.Range("H55").Value = 1329
$ws.Range("I55").Value = 539.1177
$ws.Range("J55").Value = 2549.7273
$ws.Range("K55").Value = 539.1177
$ws.Range("L55").Value = 2549.7273
$ws.Range("M55").Value = -366.1177
$ws.Range("N55").Value = -2895.7273
$ws.Range("H136").Value = 5488.2593
$ws.Range("I136").Value = 5549.3335
$ws.Range("K136").Value = 16648.0005
$ws.Range("M136").Value = -14098.0005

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 51078.91
$ws.Range("I81").Value = 74617.07000000001
$ws.Range("J81").Value = 9887.125
$ws.Range("K81").Value = 149234.14
$ws.Range("L81").Value = 19774.25
$ws.Range("M81").Value = -148173.14
$ws.Range("N81").Value = -21896.25
$ws.Range("H84").Value = 51078.91
$ws.Range("I84").Value = 74617.07000000001
$ws.Range("J84").Value = 9887.125
$ws.Range("K84").Value = 746170.7000000001
$ws.Range("L84").Value = 98871.25
$ws.Range("M84").Value = -740866.7000000001
$ws.Range("N84").Value = -109479.25
$ws.Range("H126").Value = 13990.896
$ws.Range("I126").Value = 9233.5
$ws.Range("J126").Value = 19846.154
$ws.Range("K126").Value = 27700.5
$ws.Range("L126").Value = 59538.462
$ws.Range("M126").Value = -25230.5
$ws.Range("N126").Value = -64478.462
$ws.Range("H132").Value = 2029
$ws.Range("I132").Value = 1928.0714
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5784.2142
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -3254.2142
$ws.Range("N132").Value = -12560
$ws.Range("H136").Value = 3828.125
$ws.Range("I136").Value = 2684.2354
$ws.Range("J136").Value = 4673.609
$ws.Range("K136").Value = 8052.706200000001
$ws.Range("L136").Value = 14020.827
$ws.Range("M136").Value = -5502.706200000001
$ws.Range("N136").Value = -19120.827
